# Rename six table/column headers on Sheet1 of the CCN DB Schema Revision
# workbook, matching the "Add files via upload" commit:
#   pointlog.date          -> pointlog.log_date
#   pointlog.status        -> pointlog.point_status
#   user.user_email        -> user.email
#   attendee.lastname_initial -> attendee.lastname
#   event.name              -> event.event_name
#   event.type               -> event.event_type

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I7").Value = "point_status"
$ws.Range("F7").Value = "log_date"
$ws.Range("C13").Value = "email"
$ws.Range("C18").Value = "lastname"
$ws.Range("F18").Value = "event_name"
$ws.Range("I18").Value = "event_type"

# Match the saved selection/scroll state captured in the diff.
$ws.Range("I18").Select()
$ws.Application.ActiveWindow.ScrollRow = 7
